$d = $word.ActiveDocument

$d.Content.Find.Execute("544×8=4352", $true, $true, $false, $false, $false, $true, 1, $false, "437×8=3496", 2) | Out-Null
$d.Content.Find.Execute("753×2=1506", $true, $true, $false, $false, $false, $true, 1, $false, "656×2=1312", 2) | Out-Null
$d.Content.Find.Execute("768×7=5376", $true, $true, $false, $false, $false, $true, 1, $false, "579×5=2895", 2) | Out-Null
$d.Content.Find.Execute("165×7=1155", $true, $true, $false, $false, $false, $true, 1, $false, "641×4=2564", 2) | Out-Null
$d.Content.Find.Execute("713×3=2139", $true, $true, $false, $false, $false, $true, 1, $false, "443×4=1772", 2) | Out-Null
$d.Content.Find.Execute("134×8=1072", $true, $true, $false, $false, $false, $true, 1, $false, "233×3=699", 2) | Out-Null
$d.Content.Find.Execute("364×5=1820", $true, $true, $false, $false, $false, $true, 1, $false, "339×2=678", 2) | Out-Null
$d.Content.Find.Execute("970×9=8730", $true, $true, $false, $false, $false, $true, 1, $false, "757×9=6813", 2) | Out-Null
$d.Content.Find.Execute("413×3=1239", $true, $true, $false, $false, $false, $true, 1, $false, "866×9=7794", 2) | Out-Null
$d.Content.Find.Execute("706×7=4942", $true, $true, $false, $false, $false, $true, 1, $false, "126×5=630", 2) | Out-Null
$d.Content.Find.Execute("609×9=5481", $true, $true, $false, $false, $false, $true, 1, $false, "860×2=1720", 2) | Out-Null
$d.Content.Find.Execute("815×3=2445", $true, $true, $false, $false, $false, $true, 1, $false, "374×4=1496", 2) | Out-Null
$d.Content.Find.Execute("465×3=1395", $true, $true, $false, $false, $false, $true, 1, $false, "193×3=579", 2) | Out-Null
$d.Content.Find.Execute("742×3=2226", $true, $true, $false, $false, $false, $true, 1, $false, "343×9=3087", 2) | Out-Null
$d.Content.Find.Execute("888×7=6216", $true, $true, $false, $false, $false, $true, 1, $false, "917×4=3668", 2) | Out-Null
$d.Content.Find.Execute("306×6=1836", $true, $true, $false, $false, $false, $true, 1, $false, "754×6=4524", 2) | Out-Null
$d.Content.Find.Execute("111×4=444", $true, $true, $false, $false, $false, $true, 1, $false, "486×4=1944", 2) | Out-Null
$d.Content.Find.Execute("708×6=4248", $true, $true, $false, $false, $false, $true, 1, $false, "454×9=4086", 2) | Out-Null
$d.Content.Find.Execute("472×3=1416", $true, $true, $false, $false, $false, $true, 1, $false, "933×2=1866", 2) | Out-Null
$d.Content.Find.Execute("788×2=1576", $true, $true, $false, $false, $false, $true, 1, $false, "334×6=2004", 2) | Out-Null
$d.Content.Find.Execute("175×3=525", $true, $true, $false, $false, $false, $true, 1, $false, "443×4=1772", 2) | Out-Null
$d.Content.Find.Execute("333×4=1332", $true, $true, $false, $false, $false, $true, 1, $false, "394×4=1576", 2) | Out-Null
$d.Content.Find.Execute("585×8=4680", $true, $true, $false, $false, $false, $true, 1, $false, "372×7=2604", 2) | Out-Null
$d.Content.Find.Execute("136×2=272", $true, $true, $false, $false, $false, $true, 1, $false, "156×9=1404", 2) | Out-Null
$d.Content.Find.Execute("329×3=987", $true, $true, $false, $false, $false, $true, 1, $false, "747×9=6723", 2) | Out-Null
